# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price (column D) and Volume(1h) (column E) are stored as plain text in
# this sheet. For D-column values that look like plain numbers, force the
# cell to Text format before writing so Excel doesn't auto-convert them to
# numeric values (which would silently drop formatting like trailing
# zeros, e.g. "7.150" -> 7.15). The format is reset to General right after
# so the cell's number format matches the original ("General").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.198.53'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '1.659.91'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.20'
$ws.Range('D5').NumberFormat = 'general'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5253'
$ws.Range('D6').NumberFormat = 'general'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2634'
$ws.Range('D8').NumberFormat = 'general'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06311'
$ws.Range('D9').NumberFormat = 'general'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.63'
$ws.Range('D10').NumberFormat = 'general'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07814'
$ws.Range('D11').NumberFormat = 'general'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.492'
$ws.Range('D12').NumberFormat = 'general'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').Value = '1.671.14'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').Value = '1.888.76'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5548'
$ws.Range('D15').NumberFormat = 'general'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = '0.0₅8030'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.19'
$ws.Range('D17').NumberFormat = 'general'
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = '26.218.73'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.638'
$ws.Range('D20').NumberFormat = 'general'
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '196.50'
$ws.Range('D21').NumberFormat = 'general'
$ws.Range('E21').Value = '  +1.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.13'
$ws.Range('D22').NumberFormat = 'general'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.966'
$ws.Range('D23').NumberFormat = 'general'
$ws.Range('E23').Value = '  -1.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.008'
$ws.Range('D24').NumberFormat = 'general'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.68'
$ws.Range('D25').NumberFormat = 'general'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1205'
$ws.Range('D26').NumberFormat = 'general'
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.150'
$ws.Range('D27').NumberFormat = 'general'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.05'
$ws.Range('D28').NumberFormat = 'general'
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.499'
$ws.Range('D29').NumberFormat = 'general'
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05769'
$ws.Range('D30').NumberFormat = 'general'
$ws.Range('E30').Value = '  -2.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.277'
$ws.Range('D31').NumberFormat = 'general'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.489'
$ws.Range('D32').NumberFormat = 'general'
$ws.Range('E32').Value = '  -2.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.352'
$ws.Range('D33').NumberFormat = 'general'
$ws.Range('E33').Value = '  +2.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.583'
$ws.Range('D34').NumberFormat = 'general'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9541'
$ws.Range('D35').NumberFormat = 'general'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.806'
$ws.Range('D36').NumberFormat = 'general'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5712'
$ws.Range('D38').NumberFormat = 'general'
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01598'
$ws.Range('D39').NumberFormat = 'general'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.962'
$ws.Range('D40').NumberFormat = 'general'
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('D41').Value = '1.060.66'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8509'
$ws.Range('D42').NumberFormat = 'general'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.27'
$ws.Range('D44').NumberFormat = 'general'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').Value = '1.799.71'
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.18'
$ws.Range('D46').NumberFormat = 'general'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.011'
$ws.Range('D47').NumberFormat = 'general'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.991'
$ws.Range('D49').NumberFormat = 'general'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05212'
$ws.Range('D50').NumberFormat = 'general'
$ws.Range('E50').Value = '  +0.95%  '
$ws.Range('E51').Value = '  -6.09%  '
